$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: numeric 0, bold font, centered/top aligned, thin box border
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B1").VerticalAlignment = -4160    # xlTop
$ws.Range("B1").Borders.LineStyle = 1        # xlContinuous
$ws.Range("B1").Borders.Weight = 2           # xlThin

# A2: numeric 0, same style as B1 (reuse the style via copy/paste-format
# so the generated styles.xml only gets one extra cellXfs entry)
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)          # xlPasteFormats

# B2: text label, default (unstyled) cell
$ws.Range("B2").Value = "disconnected_elements"
